$d = $word.ActiveDocument

# Set the paragraph mark's run style of the (sole, empty) paragraph to the
# "tei_supplied" character style (internal styleId "teisupplied").
$p = $d.Paragraphs.Item(1)
$p.Range.Style = $d.Styles.Item("teisupplied")
